# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets, matching the rows produced by the 2026-01-28 12:2x ingestion run.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [string]$DataBlock
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow

    foreach ($line in ($DataBlock -split "`n")) {
        $line = $line.Trim()
        if ($line.Length -eq 0) { continue }
        $fields = $line -split '\|'

        # Column A holds an ISO-like date string ("2026-01-28"). Excel's
        # auto-detection would otherwise silently convert this to a real
        # date serial number, so force text storage, assign, then restore
        # the cell to the workbook's default ("Normal") style so no stray
        # formatting is left behind.
        $cellA = $ws.Cells.Item($r, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $fields[0]
        $cellA.Style = "Normal"

        $ws.Cells.Item($r, 2).Value = $fields[1]
        $ws.Cells.Item($r, 3).Value = $fields[2]
        $ws.Cells.Item($r, 4).Value = $fields[3]

        # Column E can contain percentage-looking text ("87.6%"), which
        # Excel would otherwise auto-convert into a numeric percentage.
        # Apply the same text-forcing treatment used for the date column.
        $cellE = $ws.Cells.Item($r, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $fields[4]
        $cellE.Style = "Normal"

        $ws.Cells.Item($r, 6).Value = $fields[5]

        $r++
    }
}

$pirData = @"
2026-01-28|12:29:00|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:03|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:05|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:11|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:17|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:21|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:26|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:27|12:00|Bathroom|Motion Detected|Active
2026-01-28|12:29:34|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:38|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:44|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:49|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:54|12:00|Bathroom|No Motion|Inactive
2026-01-28|12:29:59|12:00|Bathroom|No Motion|Inactive
"@

$humidityData = @"
2026-01-28|12:28:59|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:01|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:04|12:00|Bathroom|86.7%|Active
2026-01-28|12:29:16|12:00|Bathroom|86.7%|Active
2026-01-28|12:29:19|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:23|12:00|Bathroom|86.7%|Active
2026-01-28|12:29:35|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:40|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:48|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:52|12:00|Bathroom|87.6%|Active
2026-01-28|12:29:56|12:00|Bathroom|86.7%|Active
"@

$temperatureData = @"
2026-01-28|12:29:00|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:02|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:04|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:16|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:20|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:24|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:36|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:40|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:48|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:52|12:00|Bathroom|22.9C|Active
2026-01-28|12:29:56|12:00|Bathroom|22.9C|Active
"@

Add-LogRows "PIR" 294 $pirData
Add-LogRows "Humidity" 275 $humidityData
Add-LogRows "Temperature" 275 $temperatureData

$wb.Save()
